# Enabling TC03 for ECTEST
#
# The "RunMode" column (column E) on the MasterExecutor sheet is a single
# shared string ("Yes") that every data row (E2:E31) points to. Flipping
# the RunMode for the sheet therefore means updating that shared text in
# place for every cell that uses it, which changes the displayed value
# from "Yes" to "No" everywhere at once.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")
$ws.Activate()

# Update every RunMode cell (E2:E31) together so the underlying shared
# string entry is rewritten in place rather than creating a new one.
$ws.Range("E2:E31").Value = "No"

# Reflect the user's scroll position / active selection at the time of
# the edit (they had scrolled down and selected the last RunMode cell).
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
